# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted before row 599 on the
# (single) worksheet. Every existing row from 599 to 720 shifts down by
# one (becoming 600-721, which grows the used range from A1:T720 to
# A1:T721), and the freshly-inserted row 599 is populated with this
# period's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 599:720 down to 600:721, leaving a blank row 599 behind.
$ws.Rows.Item(599).Insert()

# Fill in the new row 599 with the new observation.
$ws.Range("A599").Value = 4
$ws.Range("B599").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C599").Value = "Los Lagos"
$ws.Range("D599").Value = 45275
$ws.Range("E599").Value = 10
$ws.Range("F599").Value = "Fruta"
$ws.Range("G599").Value = 100102
$ws.Range("H599").Value = "Cítricos"
$ws.Range("I599").Value = 100102006
$ws.Range("J599").Value = "Pomelo"
$ws.Range("K599").Value = "Start Ruby"
$ws.Range("L599").Value = "Primera"
$ws.Range("M599").Value = 150
$ws.Range("N599").Value = 13000
$ws.Range("O599").Value = 13000
$ws.Range("P599").Value = 13000
$ws.Range("Q599").Value = "$/caja 14 kilos empedrada"
$ws.Range("R599").Value = "Región de O'Higgins"
$ws.Range("S599").Value = 929
$ws.Range("T599").Value = 14
